$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.001.60"
$ws.Range("E2").Value = "  -2.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.772.44"
$ws.Range("E3").Value = "  -0.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.67%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "418.41"
$ws.Range("E5").Value = "  +0.18%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.95"
$ws.Range("E6").Value = "  -8.26%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.762.57"
$ws.Range("E7").Value = "  -0.83%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  -8.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.15%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.712"
$ws.Range("E10").Value = "  -7.39%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.159"
$ws.Range("E11").Value = "  -12.50%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000341"
$ws.Range("E12").Value = "  -13.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.37"
$ws.Range("E13").Value = "  -8.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.14"
$ws.Range("E14").Value = "  +20.86%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.361.26"
$ws.Range("E15").Value = "  -1.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.00"
$ws.Range("E16").Value = "  -3.27%  "

# Row 17
$ws.Range("E17").Value = "  -1.81%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.792.71"
$ws.Range("E18").Value = "  +0.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.23"
$ws.Range("E19").Value = "  -6.79%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.091.50"
$ws.Range("E20").Value = "  -2.89%  "

# Row 21
$ws.Range("E21").Value = "  -6.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "400.80"
$ws.Range("E22").Value = "  -8.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.07"
$ws.Range("E23").Value = "  -7.75%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.57"
$ws.Range("E24").Value = "  -8.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.95"
$ws.Range("E25").Value = "  -4.54%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "36.39"
$ws.Range("E26").Value = "  -3.97%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.68"
$ws.Range("E27").Value = "  +9.96%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.14"
$ws.Range("E28").Value = "  -4.78%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.19"
$ws.Range("E29").Value = "  -6.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "695.04"
$ws.Range("E30").Value = "  -2.68%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.94"
$ws.Range("E31").Value = "  +10.32%  "

# Row 32
$ws.Range("E32").Value = "  +0.11%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.118"
$ws.Range("E33").Value = "  -3.82%  "

# Row 34
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.12"
$ws.Range("E34").Value = "  -4.23%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.148"
$ws.Range("E35").Value = "  -8.51%  "

# Row 36
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "37.03"
$ws.Range("E37").Value = "  -10.95%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "54.42"
$ws.Range("E38").Value = "  -5.88%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0750"
$ws.Range("E39").Value = "  +8.68%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0447"
$ws.Range("E40").Value = "  -8.56%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("E41").Value = "  -3.76%  "

# Row 42
$ws.Range("E42").Value = "  +0.00%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.60"
$ws.Range("E43").Value = "  +6.81%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.133"
$ws.Range("E44").Value = "  -10.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  -3.86%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.19"
$ws.Range("E46").Value = "  -3.72%  "

# Row 47
$ws.Range("E47").Value = "  -3.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("E48").Value = "  -4.39%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.07"
$ws.Range("E49").Value = "  -9.88%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.49"
$ws.Range("E50").Value = "  -4.44%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.71"
$ws.Range("E51").Value = "  -6.01%  "
